$wb = $excel.ActiveWorkbook

# --- Update the Date metadata value on the "Metadata" sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-09-02T15:43:08-05:00"

# --- Strip the "$CADSR:" prefix from the Code column on the "Concepts" sheet ---
$concepts = $wb.Worksheets.Item("Concepts")
$lastRow = $concepts.Cells.Item($concepts.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $concepts.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -ne $null -and $val.ToString().StartsWith('$CADSR:')) {
        # Force text number format so the numeric-looking code stays a text value
        # (matching the shared-string representation in the target workbook)
        # rather than being auto-converted to a numeric cell by Excel.
        $cell.NumberFormat = "@"
        $cell.Value = $val.ToString().Substring(7)
    }
}
